# Append 30 new trading days (13-11-2023 .. 22-12-2023) to the ENJSA
# historical price sheet, continuing directly after the last existing
# row (A1440 = "10-11-2023").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("isyatirim")

# Columns: Tarih, Kapanış(TL), Min(TL), Max(TL), AOF(TL), Hacim(TL),
#          Sermaye(mn TL), USDTRY, BIST 100, PiyasaDeğeri(mn TL),
#          PiyasaDeğeri(mn USD), HalkaAçık PD(mn TL), HalkaAçık PD(mn USD)
$newRows = @"
13-11-2023,49.52,49.52,52.55,51.03,670470812,1181,28.5934,7568,58487,2045,11691,409
14-11-2023,49.66,48.52,50.35,49.43,522003667,1181,28.6341,7671,58652,2048,11725,409
15-11-2023,50.6,49.12,51.95,50.43,1061519961,1181,28.6679,7666,59762,2085,11946,417
16-11-2023,50.3,50.3,51,50.55,533373271,1181,28.7045,7781,59408,2070,11876,414
17-11-2023,49.94,49.8,50.6,50.02,566275486,1181,28.7109,7853,58983,2054,11791,411
20-11-2023,50.1,49.96,50.75,50.3,516483981,1181,28.7568,7946,59172,2058,11828,411
21-11-2023,50.35,50.25,51.25,50.66,703285224,1181,28.7992,8047,59467,2065,11887,413
22-11-2023,50.05,50.05,51.05,50.42,470814371,1181,28.8136,7987,59113,2052,11817,410
23-11-2023,49.2,49.2,50.6,50.1,481201590,1181,28.8304,7899,58109,2016,11616,403
24-11-2023,49.58,48.82,49.84,49.27,418295810,1181,28.8607,7960,58557,2029,11706,406
27-11-2023,49.6,49.34,50.35,49.7,526207721,1181,28.9104,8107,58581,2026,11710,405
28-11-2023,49,48.84,49.86,49.35,494094181,1181,28.9289,8094,57872,2001,11569,400
29-11-2023,48.82,48.64,49.44,49.05,437718905,1181,28.9116,8009,57660,1994,11526,399
30-11-2023,47.56,47.56,49.04,48.41,385761544,1181,28.8728,7949,56172,1945,11229,389
01-12-2023,47.66,47.66,48.64,47.98,340720281,1181,28.9073,8026,56290,1947,11252,389
04-12-2023,48.98,47.78,48.98,48.39,470938081,1181,28.9035,8087,57849,2001,11564,400
05-12-2023,47.88,47.74,49.06,48.34,333433533,1181,28.9225,8057,56550,1955,11304,391
06-12-2023,46.32,46.22,48.08,47.2,364109587,1181,28.9313,7855,54707,1891,10936,378
07-12-2023,47.4,46.16,47.52,46.76,377294026,1181,28.926,7979,55983,1935,11191,387
08-12-2023,48.46,47.34,48.7,48.05,398144687,1181,28.9495,7914,57235,1977,11441,395
11-12-2023,47.04,46.82,48.38,47.29,288200231,1181,28.9917,7728,55557,1916,11106,383
12-12-2023,46.52,45.98,47.1,46.42,384144789,1181,29.0295,7749,54943,1893,10983,378
13-12-2023,44.6,44.6,46.9,45.37,339961866,1181,29.0579,7529,52676,1813,10530,362
14-12-2023,45.38,44.14,46,44.77,461351734,1181,29.0328,7813,53597,1846,10714,369
15-12-2023,46.7,45.42,46.7,46.3,522862682,1181,29.0369,7994,55156,1900,11026,380
18-12-2023,45.18,45.18,46.78,46.03,323744622,1181,29.0576,7790,53361,1836,10667,367
19-12-2023,44.44,43.88,45.14,44.28,358222051,1181,29.0822,7726,52487,1805,10492,361
20-12-2023,44.38,44.2,45.18,44.61,261888134,1181,29.1201,7663,52416,1800,10478,360
21-12-2023,45.26,43.82,45.36,44.43,361470010,1181,29.1422,7765,53455,1834,10686,367
22-12-2023,44.94,44.68,45.98,45.24,366072252,1181,29.1869,7558,53077,1819,10610,364
"@

$lines = $newRows.Trim() -split "`n"

$startRow = 1441
$r = $startRow
foreach ($line in $lines) {
    $parts = $line.Trim() -split ","

    # Column A is the date label, stored as text (shared string). A
    # leading apostrophe forces Excel to keep it as text instead of
    # auto-converting ambiguous dd-mm-yyyy strings into a date serial.
    $ws.Cells.Item($r, 1).Value = "'" + $parts[0]

    # Columns B..M are numeric.
    for ($c = 1; $c -lt $parts.Length; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = [double]$parts[$c]
    }

    $r = $r + 1
}
